# Update the db_rdmarcas sheet: rows 2-4 become genuine numbers (previously
# stored as text like "5001.00"), rows 5-8 get recomputed running totals
# ("Meta.AC" / "Venda.AC" / "Sobras" / "P"), and a new row 9 is appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 5001
$ws.Range("C2").Value = 5001
$ws.Range("D2").Value = 5050
$ws.Range("E2").Value = 5050
$ws.Range("F2").Value = 49
$ws.Range("G2").Value = 100.98
# Row 3
$ws.Range("B3").Value = 4000
$ws.Range("C3").Value = 9001
$ws.Range("D3").Value = 4000
$ws.Range("E3").Value = 9050
$ws.Range("F3").Value = 49
$ws.Range("G3").Value = 100.54
# Row 4
$ws.Range("B4").Value = 4000
$ws.Range("C4").Value = 13001
$ws.Range("D4").Value = 5000
$ws.Range("E4").Value = 14050
$ws.Range("F4").Value = 1049
$ws.Range("G4").Value = 108.07
# Row 5
$ws.Range("A5").Value = "'05/08/5000"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "'5000.00"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'18001.00"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'1000.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'15050.00"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'2951.00"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'83.61"
$ws.Range("G5").Style = "Normal"
# Row 6
$ws.Range("C6").Value = "'22142.00"
$ws.Range("C6").Style = "Normal"
$ws.Range("E6").Value = "'19192.00"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'2950.00"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'86.68"
$ws.Range("G6").Style = "Normal"
# Row 7
$ws.Range("C7").Value = "'26142.00"
$ws.Range("C7").Style = "Normal"
$ws.Range("E7").Value = "'23192.00"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'2950.00"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'88.72"
$ws.Range("G7").Style = "Normal"
# Row 8
$ws.Range("C8").Value = "'34142.00"
$ws.Range("C8").Style = "Normal"
$ws.Range("E8").Value = "'31192.00"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'2950.00"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'91.36"
$ws.Range("G8").Style = "Normal"
# Row 9
$ws.Range("A9").Value = "'05/08/2023"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "'8000.00"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'42142.00"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'8000.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'39192.00"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'2950.00"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'93.00"
$ws.Range("G9").Style = "Normal"
